$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns for each coin row.
# Price values are forced to text (they must not be auto-converted to
# numbers/dates by Excel), then the temporary number-format tweak is
# cleared again so the cell keeps its original (default) style.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.571.54"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -1.24%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.847.26"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  -0.88%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "314.70"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -1.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.004"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.02%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4269"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  -2.26%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3659"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.07%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "44.67"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.41%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07326"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -2.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.8858"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -5.34%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "20.88"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.72%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.861.42"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -2.64%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.362"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  -1.46%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.560"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -2.54%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06925"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +0.53%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "1.005"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.02%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "78.87"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -3.06%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008896"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.78%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "1.004"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.47"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -2.55%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "27.586.83"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -1.15%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.992"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.58%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "10.69"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -3.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.106.96"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.69%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.965"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -2.42%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "153.83"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.99"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +2.96%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "121.99"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +7.49%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "5.263"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -5.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.920"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +12.18%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.08947"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.7643"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -6.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.579"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  -5.18%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.972"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.20%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.101"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -6.58%  "
$ws.Range("E37").Value = "  -0.09%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05386"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  -2.52%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.095"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -2.40%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01950"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -1.29%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.800"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.67%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "6.923"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  -1.61%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5106"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -3.09%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1657"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -2.47%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.281"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  -5.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.06582"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.60%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4766"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.65%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "10.43"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -2.12%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.37"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -3.15%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.630"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.81%  "
